$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53: number of crashed runs
$ws.Range("F53").Value = "N. crash:"
$ws.Range("G53").Formula = '=COUNTIF(G2:G51,"True")'

# Row 54: Mean across the numeric metric columns (I:P)
$ws.Range("F54").Value = "Mean:"
$ws.Range("I54:P54").Formula = "=AVERAGE(I2:I51)"

# Row 55: Standard deviation across the numeric metric columns (I:P)
$ws.Range("F55").Value = "Standard Deviation:"
$ws.Range("I55:P55").Formula = "=STDEV.S(I2:I51)"

# Match the author's final on-screen selection/scroll position
[void]$ws.Range("H57").Select()
